$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "67.414.37"
$ws.Range("E2").Value = "  -4.89%  "
# Row 3
$ws.Range("D3").Value = "3.250.63"
$ws.Range("E3").Value = "  -8.07%  "
# Row 4
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  -0.03%  "
# Row 5
$ws.Range("D5").Value = "'590.18"
$ws.Range("E5").Value = "  -5.28%  "
# Row 6
$ws.Range("D6").Value = "'153.23"
$ws.Range("E6").Value = "  -12.74%  "
# Row 7
$ws.Range("E7").Value = "  +0.06%  "
# Row 8
$ws.Range("D8").Value = "3.242.01"
$ws.Range("E8").Value = "  -8.25%  "
# Row 9
$ws.Range("E9").Value = "  -10.76%  "
# Row 10
$ws.Range("D10").Value = "'0.172"
$ws.Range("E10").Value = "  -13.08%  "
# Row 11
$ws.Range("D11").Value = "'6.81"
$ws.Range("E11").Value = "  -5.32%  "
# Row 12
$ws.Range("D12").Value = "'0.508"
$ws.Range("E12").Value = "  -13.49%  "
# Row 13
$ws.Range("D13").Value = "'38.65"
$ws.Range("E13").Value = "  -17.41%  "
# Row 14
$ws.Range("E14").Value = "  -11.59%  "
# Row 15
$ws.Range("D15").Value = "3.768.20"
$ws.Range("E15").Value = "  -8.21%  "
# Row 16
$ws.Range("D16").Value = "67.428.84"
$ws.Range("E16").Value = "  -4.96%  "
# Row 17
$ws.Range("B17").Value = "BitcoinCash"
$ws.Range("C17").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D17").Value = "'546.73"
$ws.Range("E17").Value = "  -10.18%  "
# Row 18
$ws.Range("B18").Value = "WrappedEther"
$ws.Range("C18").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D18").Value = "3.241.92"
$ws.Range("E18").Value = "  -8.42%  "
# Row 19
$ws.Range("D19").Value = "'7.27"
$ws.Range("E19").Value = "  -13.95%  "
# Row 20
$ws.Range("D20").Value = "'0.114"
$ws.Range("E20").Value = "  -6.01%  "
# Row 21
$ws.Range("D21").Value = "'15.22"
$ws.Range("E21").Value = "  -14.65%  "
# Row 22
$ws.Range("D22").Value = "'0.768"
$ws.Range("E22").Value = "  -13.56%  "
# Row 23
$ws.Range("D23").Value = "'7.91"
$ws.Range("E23").Value = "  -12.79%  "
# Row 24
$ws.Range("D24").Value = "'85.81"
$ws.Range("E24").Value = "  -12.82%  "
# Row 25
$ws.Range("D25").Value = "'13.53"
$ws.Range("E25").Value = "  -13.99%  "
# Row 26
$ws.Range("E26").Value = "  -0.07%  "
# Row 27
$ws.Range("D27").Value = "'3.22"
$ws.Range("E27").Value = "  -15.00%  "
# Row 28
$ws.Range("E28").Value = "  -10.30%  "
# Row 29
$ws.Range("D29").Value = "'29.49"
$ws.Range("E29").Value = "  -13.28%  "
# Row 30
$ws.Range("D30").Value = "'2.14"
$ws.Range("E30").Value = "  -17.49%  "
# Row 31
$ws.Range("D31").Value = "'2.72"
$ws.Range("E31").Value = "  -11.26%  "
# Row 32
$ws.Range("E32").Value = "  -11.12%  "
# Row 33
$ws.Range("D33").Value = "'549.45"
$ws.Range("E33").Value = "  -13.18%  "
# Row 34
$ws.Range("D34").Value = "'6.63"
$ws.Range("E34").Value = "  -18.86%  "
# Row 35
$ws.Range("E35").Value = "  -16.14%  "
# Row 36
$ws.Range("E36").Value = "  -0.23%  "
# Row 37
$ws.Range("B37").Value = "VeChain"
$ws.Range("C37").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D37").Value = "'0.0447"
$ws.Range("E37").Value = "  -5.89%  "
# Row 38
$ws.Range("B38").Value = "OKB"
$ws.Range("C38").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D38").Value = "'53.85"
$ws.Range("E38").Value = "  -5.44%  "
# Row 39
$ws.Range("B39").Value = "Hedera"
$ws.Range("C39").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D39").Value = "'0.0855"
$ws.Range("E39").Value = "  -14.78%  "
# Row 40
$ws.Range("B40").Value = "Cosmos"
$ws.Range("C40").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D40").Value = "'9.24"
$ws.Range("E40").Value = "  -14.99%  "
# Row 41
$ws.Range("E41").Value = "  -11.27%  "
# Row 42
$ws.Range("D42").Value = "2.941.47"
$ws.Range("E42").Value = "  -12.69%  "
# Row 43
$ws.Range("D43").Value = "'2.63"
$ws.Range("E43").Value = "  -25.05%  "
# Row 44
$ws.Range("D44").Value = "'0.263"
$ws.Range("E44").Value = "  -16.31%  "
# Row 45
$ws.Range("D45").Value = "0.0₃0587"
$ws.Range("E45").Value = "  -20.77%  "
# Row 46
$ws.Range("D46").Value = "'26.58"
$ws.Range("E46").Value = "  -17.66%  "
# Row 47
$ws.Range("E47").Value = "  -20.05%  "
# Row 48
$ws.Range("E48").Value = "  +0.01%  "
# Row 49
$ws.Range("E49").Value = "  -17.00%  "
# Row 50
$ws.Range("D50").Value = "'126.00"
$ws.Range("E50").Value = "  -5.37%  "
# Row 51
$ws.Range("E51").Value = "  -12.66%  "
